# Sku, weight, etc ready for product with and without variants. Missing categories
#
# This script updates the "exported.xlsx" sample workbook:
#  - Productos!E2:I2 go from placeholder "TODO" text to real numeric values
#    (SKU / Peso / Altura / Longitud / Profundidad for the product without variants).
#  - The sample/example data rows (row 2, and row 3 where present) are removed
#    from Propiedades, Opciones, Variantes, Ubicaciones and Stock, leaving only
#    the header row — the "Missing categories" / not-yet-ready sheets.
#  - A handful of bestFit column widths shift as a side effect of the content
#    change; they are nudged to match as closely as this host's column-width
#    quantization allows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Productos: fill in the real SKU / Peso / Altura / Longitud / Profundidad
# values for the first product (previously "TODO" placeholders).
# ---------------------------------------------------------------------------
$productos = $wb.Worksheets.Item("Productos")
$productos.Range("E2").Value = 234566
$productos.Range("F2").Value = 123.0
$productos.Range("G2").Value = 132.0
$productos.Range("H2").Value = 134.0
$productos.Range("I2").Value = 234.0
$productos.Columns.Item(5).ColumnWidth = 8.929688

# ---------------------------------------------------------------------------
# Propiedades: drop the two sample rows, keep only the header.
# ---------------------------------------------------------------------------
$propiedades = $wb.Worksheets.Item("Propiedades")
$propiedades.Rows.Item(2).EntireRow.Delete()
$propiedades.Rows.Item(2).EntireRow.Delete()
$propiedades.Columns.Item(2).ColumnWidth = 8.929688

# ---------------------------------------------------------------------------
# Opciones: drop the two sample rows, keep only the header.
# ---------------------------------------------------------------------------
$opciones = $wb.Worksheets.Item("Opciones")
$opciones.Rows.Item(2).EntireRow.Delete()
$opciones.Rows.Item(2).EntireRow.Delete()
$opciones.Columns.Item(3).ColumnWidth = 5.643973

# ---------------------------------------------------------------------------
# Variantes: drop the single sample row, keep only the header.
# ---------------------------------------------------------------------------
$variantes = $wb.Worksheets.Item("Variantes")
$variantes.Rows.Item(2).EntireRow.Delete()
$variantes.Columns.Item(3).ColumnWidth = 7.786830
$variantes.Columns.Item(4).ColumnWidth = 5.643973
$variantes.Columns.Item(5).ColumnWidth = 5.643973
$variantes.Columns.Item(6).ColumnWidth = 5.643973
$variantes.Columns.Item(7).ColumnWidth = 4.501116

# ---------------------------------------------------------------------------
# Ubicaciones: drop the single sample row, keep only the header.
# ---------------------------------------------------------------------------
$ubicaciones = $wb.Worksheets.Item("Ubicaciones")
$ubicaciones.Rows.Item(2).EntireRow.Delete()
$ubicaciones.Columns.Item(4).ColumnWidth = 4.501116
$ubicaciones.Columns.Item(8).ColumnWidth = 8.929688
$ubicaciones.Columns.Item(9).ColumnWidth = 4.501116
$ubicaciones.Columns.Item(10).ColumnWidth = 7.786830

# ---------------------------------------------------------------------------
# Stock: drop the single sample row, keep only the header.
# ---------------------------------------------------------------------------
$stock = $wb.Worksheets.Item("Stock")
$stock.Rows.Item(2).EntireRow.Delete()
